$wb = $excel.ActiveWorkbook
$sheets = $wb.Worksheets

# ---------------------------------------------------------------------------
# 1. Insert a brand-new worksheet "2022-Q4" right after "总计" (the totals
#    sheet), pushing every quarter sheet down by one position.
# ---------------------------------------------------------------------------
$totalSheet = $sheets.Item(1)

$newSheet = $sheets.Add($null, $totalSheet)
$newSheet.Name = "2022-Q4"

# Worksheets resolve by live index, so fetch the "2022-Q3" reference sheet
# only AFTER the insert has shifted everything one slot to the right.
$refSheet = $sheets.Item("2022-Q3")

# Clone header-row (B1:H1) and column-A (A2:A8) formatting from the 2022-Q3
# sheet so the new sheet matches the workbook's established look.
$refSheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)   # xlPasteFormats

$refSheet.Range("A2:A8").Copy()
$newSheet.Range("A2:A8").PasteSpecial(-4122)   # xlPasteFormats

# Header row
$newSheet.Cells.Item(1, 2).Value = "基金代码"
$newSheet.Cells.Item(1, 3).Value = "基金名称"
$newSheet.Cells.Item(1, 4).Value = "基金规模"
$newSheet.Cells.Item(1, 5).Value = "股票总仓位"
$newSheet.Cells.Item(1, 6).Value = "仓位占比"
$newSheet.Cells.Item(1, 7).Value = "持有市值(亿元)"
$newSheet.Cells.Item(1, 8).Value = "仓位排名"

# Fund holdings for 2022-Q4 (A=index, B=code, C=name, D=scale, E=stock
# position, F=position ratio, G=holding value, H=position rank)
$fundRows = @(
    @(0, "001822", "华商智能生活灵活配置混合A",   "33.45", "90.70", "4.99", "1.6692", 5),
    @(1, "010550", "华商双擎领航混合",             "12.41", "90.98", "5.50", "0.6826", 5),
    @(2, "015385", "华商智能生活灵活配置混合C",   "11.97", "90.70", "4.99", "0.5973", 5),
    @(3, "013886", "华商新能源汽车混合A",           "9.05", "89.00", "3.97", "0.3593", 9),
    @(4, "014350", "华商卓越成长一年持有混合A",     "3.05", "93.10", "5.33", "0.1626", 4),
    @(5, "013887", "华商新能源汽车混合C",           "3.70", "89.00", "3.97", "0.1469", 9),
    @(6, "014351", "华商卓越成长一年持有混合C",     "0.10", "93.10", "5.33", "0.0053", 4)
)

for ($i = 0; $i -lt $fundRows.Length; $i++) {
    $r   = $i + 2
    $row = $fundRows[$i]

    $newSheet.Cells.Item($r, 1).Value = $row[0]
    $newSheet.Cells.Item($r, 2).Value = "'" + $row[1]   # fund code, keep as text
    $newSheet.Cells.Item($r, 3).Value = $row[2]
    $newSheet.Cells.Item($r, 4).Value = "'" + $row[3]   # scale, keep as text
    $newSheet.Cells.Item($r, 5).Value = "'" + $row[4]   # stock position, text
    $newSheet.Cells.Item($r, 6).Value = "'" + $row[5]   # position ratio, text
    $newSheet.Cells.Item($r, 7).Value = "'" + $row[6]   # holding value, text
    $newSheet.Cells.Item($r, 8).Value = $row[7]
}

# ---------------------------------------------------------------------------
# 2. Update the "总计" (totals) sheet: insert the 2022-Q4 summary row at the
#    top of the data (row 2) and shift the rest down.
# ---------------------------------------------------------------------------
$totalSheet.Cells.Item(2, 1).Copy()
$totalSheet.Range("A2:A7").PasteSpecial(-4122)   # xlPasteFormats

$summaryRows = @(
    @(0, "2022-Q4", 7, 3.62),
    @(1, "2022-Q3", 8, 4.96),
    @(2, "2022-Q2", 5, 4.09),
    @(3, "2022-Q1", 6, 3.13),
    @(4, "2021-Q4", 5, 1.43),
    @(5, "2021-Q3", 8, 1.49)
)

for ($i = 0; $i -lt $summaryRows.Length; $i++) {
    $r   = $i + 2
    $row = $summaryRows[$i]

    $totalSheet.Cells.Item($r, 1).Value = $row[0]
    $totalSheet.Cells.Item($r, 2).Value = $row[1]
    $totalSheet.Cells.Item($r, 3).Value = $row[2]
    $totalSheet.Cells.Item($r, 4).Value = $row[3]
}

# ---------------------------------------------------------------------------
# 3. Restore the originally-active tab. Inserting a sheet focuses it, but the
#    source workbook had "2021-Q3" (now the last tab) selected; put the
#    selection back so we don't introduce an unrelated view-state change.
# ---------------------------------------------------------------------------
$sheets.Item($sheets.Count).Activate()
